$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original data (before edit):
#   A1 = "Nombre del Grupo"                         B1 = "Eventos Científicos"
#   A2 = "INVEMAR - Calidad Ambiental Marina"        B2 = "1.- Seminario ... 62.- Taller ..."
#
# New filter-word rows are being added for "rede de monitoreo y calidad del agua":
# a new row is inserted above the existing INVEMAR row, and a new row is appended
# below it, so the sheet grows from 2 data rows to 4.

# Shift the existing second row ("INVEMAR...") down to row 3, unchanged.
$ws.Range("A3").Value = "INVEMAR - Calidad Ambiental Marina"
$ws.Range("B3").Value = "1.- Seminario : Seminario Internacional. El monitoreo comunitario del agua: aprendizajes socionaturales para la gestión del territorio SANTA MARTA, desde 2021-06-11 - hasta 2021-06-11 Ámbito: Internacional, Tipos de participación: Ponente Instituciones asociadas Nombre de la institución: INSTITUTO DE INVESTIGACIONES MARINAS Y COSTERAS JOSE BENITO VIVES DE ANDREIS INVEMAR  Tipo de vinculación Patrocinadora `n 62.- Taller : DISEÑO DE MONITOREO DE CALIDAD AMBIENTAL EN ZONAS MARINO-COSTERAS. Interpretación de parámetros indicadores de calidad y contaminación de aguas marinas y costeras SANTA MARTA, desde 2015-08-26 - hasta 2015-08-28 Ámbito: Nacional, Tipos de participación: Asistente , Organizador Instituciones asociadas Nombre de la institución: INSTITUTO DE INVESTIGACIONES MARINAS Y COSTERAS JOSE BENITO VIVES DE ANDREIS INVEMAR  Tipo de vinculación Patrocinadora"

# New row 2: "Ambiente y Vida" group, inserted above the INVEMAR row.
$ws.Range("A2").Value = "Ambiente y Vida"
$ws.Range("B2").Value = "217.- Taller : Curso Taller Monitoreo y Analisis de Calidad y Cantidad del agua CÚCUTA, desde 2004-01-01 - hasta Ámbito: Nacional, Tipos de participación: Organizador Instituciones asociadas Nombre de la institución: UNIVERSIDAD FRANCISCO DE PAULA SANTANDER  Tipo de vinculación Patrocinadora"

# New row 4: "Resiliencia y Saneamiento, RESA" group, appended below.
$ws.Range("A4").Value = "Resiliencia y Saneamiento, RESA"
$ws.Range("B4").Value = "27.- Taller : Taller Teorico sobre Muestreo y Monitoreo de Suelos y Aguas Subterráneas BARRANQUILLA, desde 2009-08-24 - hasta 2009-08-24 Ámbito: Nacional, Tipos de participación: Organizador Instituciones asociadas Nombre de la institución: UNIVERSIDAD NACIONAL DE COLOMBIA SEDE BOGOTA  Tipo de vinculación Patrocinadora"

# Column A carries the bold / bordered / centered-top "group name" style (same as
# the header row and the original A2 cell) - apply it to the two new group cells
# (A2, A4) plus re-apply to the shifted A3 cell so every group-name cell matches.
$a2 = $ws.Range("A2")
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$a3 = $ws.Range("A3")
$a3.Font.Bold = $true
$a3.HorizontalAlignment = -4108
$a3.VerticalAlignment = -4160
$a3.Borders.LineStyle = 1

$a4 = $ws.Range("A4")
$a4.Font.Bold = $true
$a4.HorizontalAlignment = -4108
$a4.VerticalAlignment = -4160
$a4.Borders.LineStyle = 1
